$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 842, shifting existing rows 842:883 down to 843:884
$ws.Rows.Item(842).Insert()

# The new date value looks like a date ("2026/02/19") and Excel would otherwise
# auto-convert it into a date serial number. Force the cell to text first so the
# literal string is preserved, then restore a normal (unformatted) style.
$ws.Range("A842").NumberFormat = "@"
$ws.Range("A842").Value = "2026/02/19"
$ws.Range("A842").Style = "Normal"

# Populate the remaining columns of the newly inserted row with the new data point
$ws.Range("B842").Value = "木"
$ws.Range("C842").Value = 19
$ws.Range("D842").Value = 201
